$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# (e.g. "212.23", "4.00") are not auto-converted to numbers.
$ws.Range("D2:D48").NumberFormat = "@"

$ws.Range("D2").Value = '26.067.99'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '1.596.00'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '212.23'
$ws.Range("E5").Value = '  +2.23%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.485'
$ws.Range("E7").Value = '  +1.33%  '
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").Value = '0.0616'
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '17.92'
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  +4.61%  '
$ws.Range("D12").Value = '1.817.56'
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").Value = '1.602.86'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = '4.00'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = '0.510'
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = '26.031.47'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '60.42'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D20").Value = '203.80'
$ws.Range("E20").Value = '  +9.69%  '
$ws.Range("E21").Value = '  +2.26%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  +10.33%  '
$ws.Range("D25").Value = '141.65'
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -3.27%  '
$ws.Range("D28").Value = '15.19'
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("D29").Value = '6.44'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("D31").Value = '0.0470'
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  +2.63%  '
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("D36").Value = '1.108.56'
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '0.0159'
$ws.Range("E37").Value = '  +6.44%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("D40").Value = '0.777'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("D41").Value = '0.490'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").Value = '0.777'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = '1.730.26'
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("D44").Value = '92.34'
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").Value = '1.50'
$ws.Range("E46").Value = '  +5.04%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '53.34'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0100'
$ws.Range("E48").Value = '  -9.46%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("E51").Value = '  +0.10%  '
